{"js": "const replacements = [\n  [\"853\u00d76=5118\", \"849\u00d76=5094\"],\n  [\"352\u00d79=3168\", \"682\u00d73=2046\"],\n  [\"415\u00d75=2075\", \"294\u00d77=2058\"],\n  [\"443\u00d72=886\", \"188\u00d78=1504\"],\n  [\"254\u00d79=2286\", \"370\u00d74=1480\"],\n  [\"431\u00d79=3879\", \"862\u00d79=7758\"],\n  [\"987\u00d77=6909\", \"559\u00d78=4472\"],\n  [\"238\u00d79=2142\", \"562\u00d78=4496\"],\n  [\"637\u00d78=5096\", \"582\u00d78=4656\"],\n  [\"581\u00d73=1743\", \"847\u00d73=2541\"],\n  [\"151\u00d75=755\", \"448\u00d78=3584\"],\n  [\"876\u00d73=2628\", \"938\u00d73=2814\"],\n  [\"149\u00d75=745\", \"315\u00d76=1890\"],\n  [\"131\u00d79=1179\", \"258\u00d74=1032\"],\n  [\"578\u00d78=4624\", \"766\u00d73=2298\"],\n  [\"134\u00d79=1206\", \"799\u00d79=7191\"],\n  [\"758\u00d74=3032\", \"735\u00d75=3675\"],\n  [\"875\u00d74=3500\", \"713\u00d77=4991\"],\n  [\"527\u00d79=4743\", \"591\u00d79=5319\"],\n  [\"879\u00d77=6153\", \"650\u00d73=1950\"],\n  [\"224\u00d77=1568\", \"219\u00d79=1971\"],\n  [\"196\u00d77=1372\", \"483\u00d75=2415\"],\n  [\"577\u00d74=2308\", \"218\u00d78=1744\"],\n  [\"318\u00d78=2544\", \"801\u00d78=6408\"],\n  [\"637\u00d73=1911\", \"929\u00d79=8361\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"853\u00d76=5118\", \"849\u00d76=5094\"),\n    @(\"352\u00d79=3168\", \"682\u00d73=2046\"),\n    @(\"415\u00d75=2075\", \"294\u00d77=2058\"),\n    @(\"443\u00d72=886\", \"188\u00d78=1504\"),\n    @(\"254\u00d79=2286\", \"370\u00d74=1480\"),\n    @(\"431\u00d79=3879\", \"862\u00d79=7758\"),\n    @(\"987\u00d77=6909\", \"559\u00d78=4472\"),\n    @(\"238\u00d79=2142\", \"562\u00d78=4496\"),\n    @(\"637\u00d78=5096\", \"582\u00d78=4656\"),\n    @(\"581\u00d73=1743\", \"847\u00d73=2541\"),\n    @(\"151\u00d75=755\", \"448\u00d78=3584\"),\n    @(\"876\u00d73=2628\", \"938\u00d73=2814\"),\n    @(\"149\u00d75=745\", \"315\u00d76=1890\"),\n    @(\"131\u00d79=1179\", \"258\u00d74=1032\"),\n    @(\"578\u00d78=4624\", \"766\u00d73=2298\"),\n    @(\"134\u00d79=1206\", \"799\u00d79=7191\"),\n    @(\"758\u00d74=3032\", \"735\u00d75=3675\"),\n    @(\"875\u00d74=3500\", \"713\u00d77=4991\"),\n    @(\"527\u00d79=4743\", \"591\u00d79=5319\"),\n    @(\"879\u00d77=6153\", \"650\u00d73=1950\"),\n    @(\"224\u00d77=1568\", \"219\u00d79=1971\"),\n    @(\"196\u00d77=1372\", \"483\u00d75=2415\"),\n    @(\"577\u00d74=2308\", \"218\u00d78=1744\"),\n    @(\"318\u00d78=2544\", \"801\u00d78=6408\"),\n    @(\"637\u00d73=1911\", \"929\u00d79=8361\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}"}
